$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44463
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 12000
$ws.Range("P2").Value = 1200

# Row 3
$ws.Range("D3").Value = 44425
$ws.Range("J3").Value = 30
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 13000
$ws.Range("P3").Value = 1300

# Row 4
$ws.Range("D4").Value = 44525
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 9000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 9000
$ws.Range("P4").Value = 900

# Row 5
$ws.Range("D5").Value = 44523
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 9000
$ws.Range("P5").Value = 900

# Row 7
$ws.Range("D7").Value = 44369
$ws.Range("J7").Value = 25
$ws.Range("K7").Value = 8000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 8000
$ws.Range("P7").Value = 800

# Row 8
$ws.Range("D8").Value = 44526
$ws.Range("J8").Value = 25
$ws.Range("K8").Value = 9000
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = 9000
$ws.Range("P8").Value = 900

# Row 9
$ws.Range("D9").Value = 44348
$ws.Range("J9").Value = 20
$ws.Range("K9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = 10000
$ws.Range("P9").Value = 1000

# Row 10
$ws.Range("D10").Value = 44473
$ws.Range("J10").Value = 25
$ws.Range("K10").Value = 11000
$ws.Range("L10").Value = 11000
$ws.Range("M10").Value = 11000
$ws.Range("P10").Value = 1100

# Row 11
$ws.Range("D11").Value = 44530
$ws.Range("J11").Value = 30
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 10000
$ws.Range("P11").Value = 1000
